# Adding more things for sensitivity and the final report outline
#
# 1. "Steam explosion yields": insert a new row above the old row 46,
#    add a hemicellulose-recovery line (label + formula), and move the
#    selection to the old row-46 line (now row 47, "A47" after the shift).
# 2. Normalize the lingering multi-range selections left over on the
#    "Χαρακτηρισμός της τροφοδοσίας" / "Bioreactor mass balance" /
#    "Mass Yields" sheets down to the single cell that was really the
#    point of interest.
# 3. Leave "Steam explosion yields" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Χαρακτηρισμός της τροφοδοσίας -----------------------------------
$wsFeed = $wb.Worksheets.Item(1)
$wsFeed.Range("B35").Select()

# --- Bioreactor mass balance -------------------------------------------
$wsBio = $wb.Worksheets.Item(3)
$wsBio.Range("P37").Select()

# --- Mass Yields ---------------------------------------------------------
$wsMass = $wb.Worksheets.Item(4)
$wsMass.Range("D12").Select()

# --- Steam explosion yields ----------------------------------------------
$wsSteam = $wb.Worksheets.Item(2)

# Insert a new row at 46; everything below (old 46..59) shifts to 47..60
# and formulas referencing those rows are adjusted automatically.
$wsSteam.Rows.Item(46).Insert()

# New line: hemicellulose recovery = sugars recovered (C22) / hemicellulose fed (B6)
$wsSteam.Range("A46").Value = "Ανάκτηση ημικυτταρίνης"
$wsSteam.Range("B46").Formula = "=C22/B6"

$wsSteam.Range("A47").Select()

# Make "Steam explosion yields" the active tab last, so it sticks as the
# workbook's active sheet.
$wsSteam.Activate()
